$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 6 (pushes the existing row 6 and everything
# below it down by two rows), making room for the new "StartGame" button
# entries right after the existing GUI button rows.
$ws.Range("A6:A7").EntireRow.Insert()

# Rename the existing menu-button rows to be more specific.
$ws.Range("A4").Value = "GUI_Button_Menu_Hover"
$ws.Range("A5").Value = "GUI_Button_Menu_Select"

# Fill in the two newly inserted rows with the start-game button assets.
$ws.Range("A6").Value = "GUI_Button_StartGame_Hover"
$ws.Range("A7").Value = "GUI_Button_StartGame_Select"
$ws.Range("C7").Value = "✅"
